# Generate Report for Handoff
#
# The nightly localization-status report regenerated its handoff info for
# the files that were still pending ("low" priority / not yet handed off):
# 543b4733-*, 6d6bd283-*, 6e36b6d6-*, ac4bd9b1-* now show priority "ht"
# (already handed off) with refreshed "Latest Handoff Datetime" /
# "Latest HO Xliff Generate Date" stamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# zh-cn: rows 4-7 -> Priority (E) "low" -> "ht", Latest Handoff Datetime (H)
# refreshed from 2016-08-22 02:42:34 -> 2016-08-22 02:42:51
foreach ($r in 4..7) {
    $zhcn.Range("E$r").Value = "ht"
    $zhcn.Range("H$r").Value = "2016-08-22 02:42:51"
}

# de-de: rows 4-7 -> Priority (E) "low" -> "ht", Latest Handoff Datetime (H)
# refreshed from 2016-08-22 02:42:39 -> 2016-08-22 02:42:55
foreach ($r in 4..7) {
    $dede.Range("E$r").Value = "ht"
    $dede.Range("H$r").Value = "2016-08-22 02:42:55"
}

# Overview: rows 4-7 -> Latest HO Xliff Generate Date (G) mirrors the de-de
# handoff timestamp, refreshed from 2016-08-22 02:42:39 -> 2016-08-22 02:42:55
foreach ($r in 4..7) {
    $overview.Range("G$r").Value = "2016-08-22 02:42:55"
}

Write-Host "Handoff report regenerated"
